$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.225.67"
$ws.Range("D3").Value = "3.137.11"
$ws.Range("E3").Value = "  -2.48%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'607.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'147.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.66%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.133.78"
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.36%  "
$ws.Range("E10").Value = "  -5.39%  "
$ws.Range("D11").Value = "'5.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.75%  "
$ws.Range("E12").Value = "  -5.09%  "
$ws.Range("E13").Value = "  -4.13%  "
$ws.Range("D14").Value = "'36.36"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.92%  "
$ws.Range("D15").Value = "3.648.87"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("D16").Value = "64.250.79"
$ws.Range("E16").Value = "  -3.27%  "
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "3.131.06"
$ws.Range("E18").Value = "  -2.69%  "
$ws.Range("E19").Value = "  -4.29%  "
$ws.Range("D20").Value = "'478.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.44%  "
$ws.Range("D21").Value = "'14.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.33%  "
$ws.Range("E22").Value = "  -3.20%  "
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("E24").Value = "  -5.36%  "
$ws.Range("D25").Value = "'83.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -2.93%  "
$ws.Range("E28").Value = "  -5.95%  "
$ws.Range("E29").Value = "  -5.62%  "
$ws.Range("D30").Value = "'0.123"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -20.53%  "
$ws.Range("D31").Value = "'6.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("D33").Value = "'2.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.85%  "
$ws.Range("D34").Value = "'26.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.46%  "
$ws.Range("E35").Value = "  -5.37%  "
$ws.Range("E36").Value = "  -5.37%  "
$ws.Range("D37").Value = "'54.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.65%  "
$ws.Range("E38").Value = "  +2.30%  "
$ws.Range("D39").Value = "0.0₃0725"
$ws.Range("E39").Value = "  -5.38%  "
$ws.Range("D40").Value = "'451.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.60%  "
$ws.Range("E41").Value = "  -4.48%  "
$ws.Range("E42").Value = "  -4.71%  "
$ws.Range("E43").Value = "  -3.52%  "
$ws.Range("D44").Value = "2.871.67"
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("E45").Value = "  -7.97%  "
$ws.Range("E46").Value = "  -7.47%  "
$ws.Range("D47").Value = "'26.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.69%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("E49").Value = "  -2.85%  "
$ws.Range("E50").Value = "  -2.66%  "
$ws.Range("D51").Value = "'118.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.03%  "
